$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "bjg"
$ws.Range("A23").Value = "c"
